$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that previously held a list-like string with a single real entry
# (e.g. "[-, -, 'MEC-2B-Ajustagem', -]") now collapse to a plain "-",
# removing the duplicated-teacher list representation.
$cells = @("C3","D3","C4","D4","C6","D6","C7","D7","B18","C18","D18","B19","C19","D19","B20","C20","D20","B21","C21","D21")

foreach ($addr in $cells) {
    $ws.Range($addr).Value = "-"
}
